# Insert a new price record for "Arveja Verde" (Perfection / Primera, Provincia de
# Huasco) on the Hortaliza - Mercado Mayorista Lo Valledor de Santiago sheet.
#
# The new observation belongs chronologically before the existing row 103, so a
# whole row is inserted at row 103, pushing the former rows 103-131 down to
# 104-132 (and the sheet's used range grows from A1:R131 to A1:R132).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 103:131 down to 104:132, leaving a blank row 103 to fill in.
$ws.Rows("103").Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A103").Value = 6
$ws.Range("B103").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C103").Value = "Metropolitana"
$ws.Range("D103").Value = 44463
$ws.Range("E103").Value = 13
$ws.Range("F103").Value = 100112022
$ws.Range("G103").Value = "Arveja Verde"
$ws.Range("H103").Value = "Perfection"
$ws.Range("I103").Value = "Primera"
$ws.Range("J103").Value = 290
$ws.Range("K103").Value = 26000
$ws.Range("L103").Value = 28000
$ws.Range("M103").Value = 27172
$ws.Range("N103").Value = "`$/malla 25 kilos"
$ws.Range("O103").Value = "Provincia de Huasco"
$ws.Range("P103").Value = 1087
$ws.Range("Q103").Value = 25
$ws.Range("R103").Value = "Hortaliza"
